# it_emx_deep_nesting.xlsx
#
# Commit: "all it tests in package it
#          emx files in package emx add tests: 1. self references 2. tags"
#
# For this workbook: the "deep*" packages/entities used by the IT
# (integration-test) fixture are nested one level deeper under a brand
# new top-level package called "it" - i.e. every package/entity/sheet
# name that used to start with "deep" now starts with "it_deep", and a
# new top-level package "it" (parent of "it_deep") is added as the new
# root of the hierarchy in the "packages" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Rename the six entity/sheet tabs to carry the new "it_" prefix.
#    (Tab order == file order here: tab4..tab9.)
# ---------------------------------------------------------------------
$wb.Worksheets.Item(4).Name = "it_deep_advanced_TestEntity_1"
$wb.Worksheets.Item(5).Name = "it_deep_advanced_p_TestEntity_2"
$wb.Worksheets.Item(6).Name = "it_deep_TestCategorical_1"
$wb.Worksheets.Item(7).Name = "it_deep_TestXref_1"
$wb.Worksheets.Item(8).Name = "it_deep_TestXref_2"
$wb.Worksheets.Item(9).Name = "it_deep_TestMref_1"

# ---------------------------------------------------------------------
# 2. "packages" sheet: rename existing packages and insert the new
#    top-level "it" package as the new root of the hierarchy.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

# existing rows renamed in place ...
$ws1.Range("A2").Value2 = "it_deep_advanced_p"
$ws1.Range("C2").Value2 = "it_deep_advanced"

$ws1.Range("A3").Value2 = "it_deep_advanced"
$ws1.Range("B3").Value2 = ""
$ws1.Range("C3").Value2 = "it_deep"

$ws1.Range("A4").Value2 = "it_deep"
$ws1.Range("C4").Value2 = "it"

# ... plus a brand new row 5 for the new top-level "it" package
$ws1.Range("A5").Value2 = "it"

# ---------------------------------------------------------------------
# 3. "entities" sheet: every entity's package gets the new "it_" prefix.
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value2 = "it_deep"
$ws2.Range("B3").Value2 = "it_deep_advanced"
$ws2.Range("E3").Value2 = "it_deep_TestEntity_0"
$ws2.Range("B4").Value2 = "it_deep_advanced_p"
$ws2.Range("E4").Value2 = "it_deep_advanced_TestEntity_1"
$ws2.Range("B5").Value2 = "it_deep"
$ws2.Range("B6").Value2 = "it_deep"
$ws2.Range("B7").Value2 = "it_deep"
$ws2.Range("B8").Value2 = "it_deep"

# ---------------------------------------------------------------------
# 4. "attributes" sheet: every referenced entity name gets "it_" too.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value2 = "it_deep_TestEntity_0"
$ws3.Range("B3").Value2 = "it_deep_TestEntity_0"
$ws3.Range("B4").Value2 = "it_deep_TestEntity_0"
$ws3.Range("B5").Value2 = "it_deep_TestEntity_0"
$ws3.Range("D5").Value2 = "it_deep_TestXref_1"
$ws3.Range("B6").Value2 = "it_deep_advanced_TestEntity_1"
$ws3.Range("D6").Value2 = "it_deep_TestMref_1"
$ws3.Range("B7").Value2 = "it_deep_advanced_TestEntity_1"
$ws3.Range("B8").Value2 = "it_deep_advanced_p_TestEntity_2"
$ws3.Range("D8").Value2 = "it_deep_TestCategorical_1"
$ws3.Range("B9").Value2 = "it_deep_advanced_p_TestEntity_2"
$ws3.Range("B10").Value2 = "it_deep_advanced_p_TestEntity_2"
$ws3.Range("D10").Value2 = "it_deep_TestXref_1"
$ws3.Range("B11").Value2 = "it_deep_advanced_p_TestEntity_2"
$ws3.Range("B12").Value2 = "it_deep_TestCategorical_1"
$ws3.Range("B13").Value2 = "it_deep_TestCategorical_1"
$ws3.Range("B14").Value2 = "it_deep_TestXref_1"
$ws3.Range("B15").Value2 = "it_deep_TestXref_1"
$ws3.Range("B16").Value2 = "it_deep_TestXref_1"
$ws3.Range("D16").Value2 = "it_deep_TestXref_2"
$ws3.Range("B17").Value2 = "it_deep_TestXref_2"
$ws3.Range("B18").Value2 = "it_deep_TestXref_2"
$ws3.Range("D18").Value2 = "it_deep_TestMref_1"
$ws3.Range("B19").Value2 = "it_deep_TestXref_2"
$ws3.Range("B20").Value2 = "it_deep_TestMref_1"
$ws3.Range("B21").Value2 = "it_deep_TestMref_1"
$ws3.Range("D21").Value2 = "it_deep_TestCategorical_1"

# ---------------------------------------------------------------------
# 5. Column width touch-ups (bestFit-style autosize growth caused by
#    the now-longer "it_deep_..." strings).
# ---------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 15.333333333333334
$ws2.Columns.Item(5).ColumnWidth = 26.5
$ws3.Columns.Item(4).ColumnWidth = 22.166666666666668

# ---------------------------------------------------------------------
# 6. Re-create the recorded cell selections on every sheet that moved.
#    Whichever sheet we select last ends up "active" (tabSelected) -
#    re-activate "packages" at the very end to match the original tab.
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Activate()
$ws4.Range("D33").Select()

$ws7 = $wb.Worksheets.Item(7)
$ws7.Activate()
$ws7.Range("K33").Select()

$ws8 = $wb.Worksheets.Item(8)
$ws8.Activate()
$ws8.Range("K36").Select()

$ws3.Activate()
$ws3.Range("B3").Select()

$ws2.Activate()
$ws2.Range("B5").Select()

$ws1.Activate()
$ws1.Range("A2").Select()
